$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cryptos price/volume table (columns D and E, rows 2-51).
# For D-column values that look like plain numbers (e.g. "7.70", "0.0000114"),
# Excel's normal Value assignment would silently reinterpret them as floating
# point numbers and lose formatting (trailing zeros, decimal grouping, etc.).
# To preserve the exact text as it appears in the source data, we temporarily
# force those cells to Text format before assigning the value, then restore
# the default "Normal" style so no stray formatting is left behind.

$ws.Range("D2").Value = "62.124.06"
$ws.Range("E2").Value = "  -1.30%  "
$ws.Range("D3").Value = "3.425.48"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.70%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +1.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "8.08"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.43%  "
$ws.Range("E10").Value = "  +1.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.419"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.04%  "
$ws.Range("D12").Value = "4.012.65"
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.66%  "
$ws.Range("E15").Value = "  +1.03%  "
$ws.Range("D16").Value = "3.428.04"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").Value = "62.117.09"
$ws.Range("E17").Value = "  -1.39%  "
$ws.Range("E18").Value = "  +3.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "383.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("E22").Value = "  +2.50%  "
$ws.Range("E23").Value = "  +2.25%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").Value = "3.561.65"
$ws.Range("E25").Value = "  -0.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000114"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.177"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.49%  "
$ws.Range("E30").Value = "  +1.09%  "
$ws.Range("E31").Value = "  -2.20%  "
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("E33").Value = "  +0.34%  "
$ws.Range("E34").Value = "  +1.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.56"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.36%  "
$ws.Range("E36").Value = "  +1.50%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.99"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "168.71"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "30.98"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.71%  "
$ws.Range("D40").Value = "3.462.18"
$ws.Range("E40").Value = "  -0.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0784"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.79"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.782"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.79%  "
$ws.Range("E45").Value = "  -1.83%  "
$ws.Range("E46").Value = "  -1.90%  "
$ws.Range("D47").Value = "2.557.60"
$ws.Range("E47").Value = "  -0.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.51"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.78%  "
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.06%  "
$ws.Range("E51").Value = "  +0.00%  "
